# Generate Report for Handoff
#
# The "d237e0b9-97fe-4e5c-bdc2-98edb19819f6" file finished a new handoff
# cycle, so its row moves from the 2nd data row to the bottom (4th) data
# row on every sheet, its Status changes from "Handed back: in sync with
# en-US" to "Ready for handoff", and its "Latest Handoff Date(time)"
# is refreshed. Rows for the other two files shift up to fill the gap.
# Hyperlink targets stay pinned to their original cell position (the
# underlying relationship Ids are not remapped), only the visible
# "display" text of each hyperlink is refreshed to match the new cell
# text.

$wb = $excel.ActiveWorkbook

function Set-RowValues {
    param($ws, $row, $values)
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}

function Rebuild-Hyperlinks {
    param($ws, $links)
    # $links is an ordered list of hashtables: @{ Ref=...; Url=...; Display=... }
    $ws.Cells.Hyperlinks.Delete()
    foreach ($link in $links) {
        [void]$ws.Hyperlinks.Add($ws.Range($link.Ref), $link.Url, "", "", $link.Display)
    }
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-RowValues $wsOverview 2 @{
    A = "ffff0dc19648-959a-4389-8364-2c5b86ef7bed.md"
    B = "Handed back: in sync with en-US"
    C = "Handed back: in sync with en-US"
    D = "2016-03-24 05:10:50"
}
Set-RowValues $wsOverview 3 @{
    A = "ffffff2376497f-bb4a-447c-871a-202617a8dc95.md"
    B = "Handed back: in sync with en-US"
    C = "Handed back: in sync with en-US"
    D = "2016-03-24 05:10:50"
}
Set-RowValues $wsOverview 4 @{
    A = "d237e0b9-97fe-4e5c-bdc2-98edb19819f6.md"
    B = "Ready for handoff"
    C = "Ready for handoff"
    D = "2016-03-24 05:13:51"
}

Rebuild-Hyperlinks $wsOverview @(
    @{ Ref = "A2"; Url = "https://github.com/OpenLocalizationTest/oltest/blob/a48cf8f0c006f1a400f6664bd5a26003bb9c1748/e2e/d237e0b9-97fe-4e5c-bdc2-98edb19819f6.md"; Display = "ffff0dc19648-959a-4389-8364-2c5b86ef7bed.md" }
    @{ Ref = "A3"; Url = "https://github.com/OpenLocalizationTest/oltest/blob/a48cf8f0c006f1a400f6664bd5a26003bb9c1748/e2e/ffff0dc19648-959a-4389-8364-2c5b86ef7bed.md"; Display = "ffffff2376497f-bb4a-447c-871a-202617a8dc95.md" }
    @{ Ref = "A4"; Url = "https://github.com/OpenLocalizationTest/oltest/blob/a48cf8f0c006f1a400f6664bd5a26003bb9c1748/e2e/ffffff2376497f-bb4a-447c-871a-202617a8dc95.md"; Display = "d237e0b9-97fe-4e5c-bdc2-98edb19819f6.md" }
)

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

Set-RowValues $wsZh 2 @{
    A = "ffff0dc19648-959a-4389-8364-2c5b86ef7bed.md"
    B = ".md"
    C = "Handed back: in sync with en-US"
    D = "aac1a086-635d-4bc7-8d61-23c24bef0bb6.c31798f111cbfeadbf8bd1000277a568068217a3.zh-cn.xlf"
    E = "2016-03-24 05:10:46"
    F = "aac1a086-635d-4bc7-8d61-23c24bef0bb6.md"
    G = "aac1a086-635d-4bc7-8d61-23c24bef0bb6.c31798f111cbfeadbf8bd1000277a568068217a3.zh-cn.xlf"
    H = "2016-03-24 05:11:09"
    J = "Include"
}
Set-RowValues $wsZh 3 @{
    A = "ffffff2376497f-bb4a-447c-871a-202617a8dc95.md"
    B = ".md"
    C = "Handed back: in sync with en-US"
    D = "aac1a086-635d-4bc7-8d61-23c24bef0bb6.c31798f111cbfeadbf8bd1000277a568068217a3.zh-cn.xlf"
    E = "2016-03-24 05:10:46"
    F = "aac1a086-635d-4bc7-8d61-23c24bef0bb6.md"
    G = "aac1a086-635d-4bc7-8d61-23c24bef0bb6.c31798f111cbfeadbf8bd1000277a568068217a3.zh-cn.xlf"
    H = "2016-03-24 05:11:09"
    J = "Include"
}
Set-RowValues $wsZh 4 @{
    A = "d237e0b9-97fe-4e5c-bdc2-98edb19819f6.md"
    B = ".md"
    C = "Ready for handoff"
    D = "d237e0b9-97fe-4e5c-bdc2-98edb19819f6.429eb12fb805ee3757f85338cdff9edb747a06cb.zh-cn.xlf"
    E = "2016-03-24 05:13:46"
    F = "d237e0b9-97fe-4e5c-bdc2-98edb19819f6.md"
    G = "d237e0b9-97fe-4e5c-bdc2-98edb19819f6.429eb12fb805ee3757f85338cdff9edb747a06cb.zh-cn.xlf"
    H = "2016-03-24 05:12:54"
    J = "Include"
}

Rebuild-Hyperlinks $wsZh @(
    @{ Ref = "A2"; Url = "https://github.com/OpenLocalizationTest/oltest/blob/a48cf8f0c006f1a400f6664bd5a26003bb9c1748/e2e/d237e0b9-97fe-4e5c-bdc2-98edb19819f6.md"; Display = "ffff0dc19648-959a-4389-8364-2c5b86ef7bed.md" }
    @{ Ref = "D2"; Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a07925b3b42939fe8e8ce1b11d3b1ce8d17ed799/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d237e0b9-97fe-4e5c-bdc2-98edb19819f6.429eb12fb805ee3757f85338cdff9edb747a06cb.zh-cn.xlf"; Display = "aac1a086-635d-4bc7-8d61-23c24bef0bb6.c31798f111cbfeadbf8bd1000277a568068217a3.zh-cn.xlf" }
    @{ Ref = "F2"; Url = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/24442f97aee2a2eb5655c943f8fca535801a377c/e2e/d237e0b9-97fe-4e5c-bdc2-98edb19819f6.md"; Display = "aac1a086-635d-4bc7-8d61-23c24bef0bb6.md" }
    @{ Ref = "G2"; Url = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6e381c6c656d9c9eb9d31f4487840a1388418951/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d237e0b9-97fe-4e5c-bdc2-98edb19819f6.429eb12fb805ee3757f85338cdff9edb747a06cb.zh-cn.xlf"; Display = "aac1a086-635d-4bc7-8d61-23c24bef0bb6.c31798f111cbfeadbf8bd1000277a568068217a3.zh-cn.xlf" }

    @{ Ref = "A3"; Url = "https://github.com/OpenLocalizationTest/oltest/blob/a48cf8f0c006f1a400f6664bd5a26003bb9c1748/e2e/ffff0dc19648-959a-4389-8364-2c5b86ef7bed.md"; Display = "ffffff2376497f-bb4a-447c-871a-202617a8dc95.md" }
    @{ Ref = "D3"; Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5c2d4dd2aa8b9ac91fc86be8c71a26198b1d41fb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/aac1a086-635d-4bc7-8d61-23c24bef0bb6.c31798f111cbfeadbf8bd1000277a568068217a3.zh-cn.xlf"; Display = "aac1a086-635d-4bc7-8d61-23c24bef0bb6.c31798f111cbfeadbf8bd1000277a568068217a3.zh-cn.xlf" }
    @{ Ref = "F3"; Url = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fff76c77fcaf91c26042f533e30acdf1d895b243/e2e/aac1a086-635d-4bc7-8d61-23c24bef0bb6.md"; Display = "aac1a086-635d-4bc7-8d61-23c24bef0bb6.md" }
    @{ Ref = "G3"; Url = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c770b44521a56433b64c0bcb508c97a86f8435c4/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/aac1a086-635d-4bc7-8d61-23c24bef0bb6.c31798f111cbfeadbf8bd1000277a568068217a3.zh-cn.xlf"; Display = "aac1a086-635d-4bc7-8d61-23c24bef0bb6.c31798f111cbfeadbf8bd1000277a568068217a3.zh-cn.xlf" }

    @{ Ref = "A4"; Url = "https://github.com/OpenLocalizationTest/oltest/blob/a48cf8f0c006f1a400f6664bd5a26003bb9c1748/e2e/ffffff2376497f-bb4a-447c-871a-202617a8dc95.md"; Display = "d237e0b9-97fe-4e5c-bdc2-98edb19819f6.md" }
    @{ Ref = "D4"; Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5c2d4dd2aa8b9ac91fc86be8c71a26198b1d41fb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/aac1a086-635d-4bc7-8d61-23c24bef0bb6.c31798f111cbfeadbf8bd1000277a568068217a3.zh-cn.xlf"; Display = "d237e0b9-97fe-4e5c-bdc2-98edb19819f6.429eb12fb805ee3757f85338cdff9edb747a06cb.zh-cn.xlf" }
    @{ Ref = "F4"; Url = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fff76c77fcaf91c26042f533e30acdf1d895b243/e2e/aac1a086-635d-4bc7-8d61-23c24bef0bb6.md"; Display = "d237e0b9-97fe-4e5c-bdc2-98edb19819f6.md" }
    @{ Ref = "G4"; Url = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c770b44521a56433b64c0bcb508c97a86f8435c4/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/aac1a086-635d-4bc7-8d61-23c24bef0bb6.c31798f111cbfeadbf8bd1000277a568068217a3.zh-cn.xlf"; Display = "d237e0b9-97fe-4e5c-bdc2-98edb19819f6.429eb12fb805ee3757f85338cdff9edb747a06cb.zh-cn.xlf" }
)

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

Set-RowValues $wsDe 2 @{
    A = "ffff0dc19648-959a-4389-8364-2c5b86ef7bed.md"
    B = ".md"
    C = "Handed back: in sync with en-US"
    D = "aac1a086-635d-4bc7-8d61-23c24bef0bb6.c31798f111cbfeadbf8bd1000277a568068217a3.de-de.xlf"
    E = "2016-03-24 05:10:50"
    F = "aac1a086-635d-4bc7-8d61-23c24bef0bb6.md"
    G = "aac1a086-635d-4bc7-8d61-23c24bef0bb6.c31798f111cbfeadbf8bd1000277a568068217a3.de-de.xlf"
    H = "2016-03-24 05:11:16"
    J = "Include"
}
Set-RowValues $wsDe 3 @{
    A = "ffffff2376497f-bb4a-447c-871a-202617a8dc95.md"
    B = ".md"
    C = "Handed back: in sync with en-US"
    D = "aac1a086-635d-4bc7-8d61-23c24bef0bb6.c31798f111cbfeadbf8bd1000277a568068217a3.de-de.xlf"
    E = "2016-03-24 05:10:50"
    F = "aac1a086-635d-4bc7-8d61-23c24bef0bb6.md"
    G = "aac1a086-635d-4bc7-8d61-23c24bef0bb6.c31798f111cbfeadbf8bd1000277a568068217a3.de-de.xlf"
    H = "2016-03-24 05:11:16"
    J = "Include"
}
Set-RowValues $wsDe 4 @{
    A = "d237e0b9-97fe-4e5c-bdc2-98edb19819f6.md"
    B = ".md"
    C = "Ready for handoff"
    D = "d237e0b9-97fe-4e5c-bdc2-98edb19819f6.429eb12fb805ee3757f85338cdff9edb747a06cb.de-de.xlf"
    E = "2016-03-24 05:13:51"
    F = "d237e0b9-97fe-4e5c-bdc2-98edb19819f6.md"
    G = "d237e0b9-97fe-4e5c-bdc2-98edb19819f6.429eb12fb805ee3757f85338cdff9edb747a06cb.de-de.xlf"
    H = "2016-03-24 05:13:00"
    J = "Include"
}

Rebuild-Hyperlinks $wsDe @(
    @{ Ref = "A2"; Url = "https://github.com/OpenLocalizationTest/oltest/blob/a48cf8f0c006f1a400f6664bd5a26003bb9c1748/e2e/d237e0b9-97fe-4e5c-bdc2-98edb19819f6.md"; Display = "ffff0dc19648-959a-4389-8364-2c5b86ef7bed.md" }
    @{ Ref = "D2"; Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/89a197ce6185a00ee70edb94248f81f93310e673/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d237e0b9-97fe-4e5c-bdc2-98edb19819f6.429eb12fb805ee3757f85338cdff9edb747a06cb.de-de.xlf"; Display = "aac1a086-635d-4bc7-8d61-23c24bef0bb6.c31798f111cbfeadbf8bd1000277a568068217a3.de-de.xlf" }
    @{ Ref = "F2"; Url = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/933d08e209cbebbe7b0372f92d4e9aff4eb4799a/e2e/d237e0b9-97fe-4e5c-bdc2-98edb19819f6.md"; Display = "aac1a086-635d-4bc7-8d61-23c24bef0bb6.md" }
    @{ Ref = "G2"; Url = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f77ab12d7558c4afaa9e4a3d0c15025a4cf060ae/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d237e0b9-97fe-4e5c-bdc2-98edb19819f6.429eb12fb805ee3757f85338cdff9edb747a06cb.de-de.xlf"; Display = "aac1a086-635d-4bc7-8d61-23c24bef0bb6.c31798f111cbfeadbf8bd1000277a568068217a3.de-de.xlf" }

    @{ Ref = "A3"; Url = "https://github.com/OpenLocalizationTest/oltest/blob/a48cf8f0c006f1a400f6664bd5a26003bb9c1748/e2e/ffff0dc19648-959a-4389-8364-2c5b86ef7bed.md"; Display = "ffffff2376497f-bb4a-447c-871a-202617a8dc95.md" }
    @{ Ref = "D3"; Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c37eaa28eece4e2c4b49aac7f360ad74cd203c87/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/aac1a086-635d-4bc7-8d61-23c24bef0bb6.c31798f111cbfeadbf8bd1000277a568068217a3.de-de.xlf"; Display = "aac1a086-635d-4bc7-8d61-23c24bef0bb6.c31798f111cbfeadbf8bd1000277a568068217a3.de-de.xlf" }
    @{ Ref = "F3"; Url = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/6f658ac284872c856a471b01f65d920f7f678912/e2e/aac1a086-635d-4bc7-8d61-23c24bef0bb6.md"; Display = "aac1a086-635d-4bc7-8d61-23c24bef0bb6.md" }
    @{ Ref = "G3"; Url = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/caa702a4f97bc0ef13cd156b5dcbe082bf959c70/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/aac1a086-635d-4bc7-8d61-23c24bef0bb6.c31798f111cbfeadbf8bd1000277a568068217a3.de-de.xlf"; Display = "aac1a086-635d-4bc7-8d61-23c24bef0bb6.c31798f111cbfeadbf8bd1000277a568068217a3.de-de.xlf" }

    @{ Ref = "A4"; Url = "https://github.com/OpenLocalizationTest/oltest/blob/a48cf8f0c006f1a400f6664bd5a26003bb9c1748/e2e/ffffff2376497f-bb4a-447c-871a-202617a8dc95.md"; Display = "d237e0b9-97fe-4e5c-bdc2-98edb19819f6.md" }
    @{ Ref = "D4"; Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c37eaa28eece4e2c4b49aac7f360ad74cd203c87/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/aac1a086-635d-4bc7-8d61-23c24bef0bb6.c31798f111cbfeadbf8bd1000277a568068217a3.de-de.xlf"; Display = "d237e0b9-97fe-4e5c-bdc2-98edb19819f6.429eb12fb805ee3757f85338cdff9edb747a06cb.de-de.xlf" }
    @{ Ref = "F4"; Url = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/6f658ac284872c856a471b01f65d920f7f678912/e2e/aac1a086-635d-4bc7-8d61-23c24bef0bb6.md"; Display = "d237e0b9-97fe-4e5c-bdc2-98edb19819f6.md" }
    @{ Ref = "G4"; Url = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/caa702a4f97bc0ef13cd156b5dcbe082bf959c70/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/aac1a086-635d-4bc7-8d61-23c24bef0bb6.c31798f111cbfeadbf8bd1000277a568068217a3.de-de.xlf"; Display = "d237e0b9-97fe-4e5c-bdc2-98edb19819f6.429eb12fb805ee3757f85338cdff9edb747a06cb.de-de.xlf" }
)

Write-Host "Done applying handoff report update."
